# Fix terminology in cache eviction evaluation (slide 36).
#
# "K'/V'" (evicted/dirty cache line) is renamed to "T/S" (tag/state),
# matching the "K/V" -> "T/S" rename already used for the resident line.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(36)

# "Evict K’,V’cache"  ->  "Evict T,Scache"   (TextBox 73)
$evict = $s.Shapes.Item(6)
$evictRange = $evict.TextFrame.TextRange
$evictRange.Characters(7, 3).Text = ""                  # remove the trailing "K',"
$evict.TextFrame.TextRange.Characters(7, 2).Text = "T,S" # "V'" -> "T,S"
$evict.Height = 41.19843                                 # undo autofit re-wrap

# "K’"  ->  "T"   (TextBox 37)
$kprime = $s.Shapes.Item(7)
$kprime.TextFrame.TextRange.Text = "T"
$kprime.Width = 31.83295

# "V’back"  ->  "Sback"   (TextBox 38)
$vprime = $s.Shapes.Item(8)
$vprime.TextFrame.TextRange.Characters(1, 2).Text = "S"
$vprime.Width = 73.233335

# "K"  ->  "T"   (TextBox 41)
$k = $s.Shapes.Item(9)
$k.TextFrame.TextRange.Text = "T"
$k.Width = 31.83295

# "V0"  ->  "S0"   (TextBox 44)
$v = $s.Shapes.Item(10)
$v.TextFrame.TextRange.Characters(1, 1).Text = "S"
